$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9298247212236502,
    0.7246052329558914,
    0.9959094076139086,
    0.8681332129985964,
    0.9164529386717214,
    0.469262412304898,
    1.841566075300257,
    0.01774961222791054,
    1.059647248919488,
    0.5386984305736787,
    0.3376508542801637,
    0.6850273077074358,
    1.129554360817876,
    0.7141903383985572,
    75.51318630499273,
    120.6115918251162
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
